$d = $word.ActiveDocument

# 1. Insert a new bullet "What is the purpose of this analysis?" right after the
#    first bullet ("Things to do in the report." is paragraph 1; the bullet list
#    begins at paragraph 2). Inserting after paragraph 1 and giving the new
#    paragraph the same list level (top level, ilvl=0 / ListLevelNumber=1) as the
#    paragraph that used to open the list reproduces the bullet's formatting.
$firstPara = $d.Paragraphs.Item(1)
$insertRange = $firstPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newPara.Style = "List Paragraph"
$newPara.Range.ListFormat.ListLevelNumber = 1
$newPara.Range.Text = "What is the purpose of this analysis?"

# 2. Move the hidden "_GoBack" bookmark from the end of the document (after
#    "Show results of simulation") to just before the run of text in what is now
#    the third paragraph ("Describe the relevant aspects of football").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$targetPara = $d.Paragraphs.Item(3)
$bmRange = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Demote "Football Power Index (ESPN)" one level, from ilvl=1 to ilvl=2
#    (ListLevelNumber 2 -> 3), to match its siblings "Company"/"Creation date".
$found = $d.Content
$found.Find.Execute("Football Power Index (ESPN)", $true, $false, $false,
                     $false, $false, $true, 1, $false, "", 0)
$found.ListFormat.ListLevelNumber = 3
